# Auto-generated edit script: updates market-price derived cells
# across all 8 sheets per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3025.549
$ws.Range("I76").Value = 3004.5112
$ws.Range("K76").Value = 3004.5112
$ws.Range("M76").Value = -2689.5112
$ws.Range("H79").Value = 3025.549
$ws.Range("I79").Value = 3004.5112
$ws.Range("K79").Value = 3004.5112
$ws.Range("M79").Value = -1912.5112
$ws.Range("H132").Value = 2463.3845
$ws.Range("I132").Value = 2088
$ws.Range("J132").Value = 4040
$ws.Range("K132").Value = 6264
$ws.Range("L132").Value = 12120
$ws.Range("M132").Value = -3734
$ws.Range("N132").Value = -17180
$ws.Range("H137").Value = 3202.4285
$ws.Range("I137").Value = 2883.4
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 8650.200000000001
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -6100.200000000001
$ws.Range("N137").Value = -17100

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 15749.75
$ws.Range("J52").Value = 15749.75
$ws.Range("L52").Value = 15749.75
$ws.Range("N52").Value = -16385.75
$ws.Range("H112").Value = 20567
$ws.Range("J112").Value = 20567
$ws.Range("L112").Value = 20567
$ws.Range("N112").Value = -23521
$ws.Range("H137").Value = 44983.168
$ws.Range("J137").Value = 46179.8
$ws.Range("L137").Value = 46179.8
$ws.Range("N137").Value = -56379.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 29550
$ws.Range("J51").Value = 29550
$ws.Range("L51").Value = 29550
$ws.Range("N51").Value = -30532
$ws.Range("H55").Value = 29726.666
$ws.Range("J55").Value = 29726.666
$ws.Range("L55").Value = 29726.666
$ws.Range("N55").Value = -30272.666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 352.2143
$ws.Range("I22").Value = 275.66666
$ws.Range("J22").Value = 490
$ws.Range("K22").Value = 275.66666
$ws.Range("L22").Value = 490
$ws.Range("M22").Value = 74.33334000000002
$ws.Range("N22").Value = -1190
$ws.Range("H58").Value = 1421.2307
$ws.Range("I58").Value = 1357.52
$ws.Range("J58").Value = 3014
$ws.Range("K58").Value = 1357.52
$ws.Range("L58").Value = 3014
$ws.Range("M58").Value = -1154.52
$ws.Range("N58").Value = -3420
$ws.Range("H136").Value = 1421.2307
$ws.Range("I136").Value = 1357.52
$ws.Range("J136").Value = 3014
$ws.Range("K136").Value = 4072.56
$ws.Range("L136").Value = 9042
$ws.Range("M136").Value = -1522.56
$ws.Range("N136").Value = -14142

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 3866.6667
$ws.Range("J100").Value = 3866.6667
$ws.Range("L100").Value = 11600.0001
$ws.Range("N100").Value = -13222.0001
$ws.Range("H109").Value = 2197.6875
$ws.Range("I109").Value = 944.5
$ws.Range("J109").Value = 2949.6
$ws.Range("K109").Value = 2833.5
$ws.Range("L109").Value = 8848.799999999999
$ws.Range("M109").Value = -1793.5
$ws.Range("N109").Value = -10928.8
$ws.Range("H115").Value = 2859
$ws.Range("I115").Value = 1008.6667
$ws.Range("J115").Value = 3363.6365
$ws.Range("K115").Value = 3026.0001
$ws.Range("L115").Value = 10090.9095
$ws.Range("M115").Value = -1851.0001
$ws.Range("N115").Value = -12440.9095
$ws.Range("H124").Value = 2287.7778
$ws.Range("I124").Value = 930
$ws.Range("J124").Value = 2966.6667
$ws.Range("K124").Value = 2790
$ws.Range("L124").Value = 8900.000100000001
$ws.Range("M124").Value = 2120
$ws.Range("N124").Value = -18720.0001
$ws.Range("H125").Value = 3202.2693
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 3290.36
$ws.Range("K125").Value = 3000
$ws.Range("L125").Value = 9871.08
$ws.Range("M125").Value = 1920
$ws.Range("N125").Value = -19711.08
$ws.Range("H126").Value = 3485.8333
$ws.Range("I126").Value = 1866
$ws.Range("J126").Value = 4642.857
$ws.Range("K126").Value = 5598
$ws.Range("L126").Value = 13928.571
$ws.Range("M126").Value = -658
$ws.Range("N126").Value = -23808.571
$ws.Range("H130").Value = 3478.9
$ws.Range("I130").Value = 1394.5
$ws.Range("J130").Value = 4000
$ws.Range("K130").Value = 4183.5
$ws.Range("L130").Value = 12000
$ws.Range("M130").Value = 836.5
$ws.Range("N130").Value = -22040
$ws.Range("H131").Value = 845.28
$ws.Range("I131").Value = 424.0909
$ws.Range("J131").Value = 897.3371
$ws.Range("K131").Value = 1272.2727
$ws.Range("L131").Value = 2692.0113
$ws.Range("M131").Value = 3767.7273
$ws.Range("N131").Value = -12772.0113

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 22300
$ws.Range("J51").Value = 22300
$ws.Range("L51").Value = 22300
$ws.Range("N51").Value = -23318
$ws.Range("H63").Value = 26000
$ws.Range("J63").Value = 26000
$ws.Range("L63").Value = 26000
$ws.Range("N63").Value = -27372
$ws.Range("H66").Value = 26000
$ws.Range("J66").Value = 26000
$ws.Range("L66").Value = 78000
$ws.Range("N66").Value = -84864
$ws.Range("H82").Value = 30000
$ws.Range("J82").Value = 30000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30766
$ws.Range("H85").Value = 30000
$ws.Range("J85").Value = 30000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32652
$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 40786.04
$ws.Range("I22").Value = 167231.67
$ws.Range("J22").Value = 855.8421
$ws.Range("K22").Value = 167231.67
$ws.Range("L22").Value = 855.8421
$ws.Range("M22").Value = -166936.67
$ws.Range("N22").Value = -1445.8421
$ws.Range("H27").Value = 40786.04
$ws.Range("I27").Value = 167231.67
$ws.Range("J27").Value = 855.8421
$ws.Range("K27").Value = 167231.67
$ws.Range("L27").Value = 855.8421
$ws.Range("M27").Value = -167124.67
$ws.Range("N27").Value = -1069.8421
$ws.Range("H64").Value = 30750
$ws.Range("J64").Value = 30750
$ws.Range("L64").Value = 30750
$ws.Range("N64").Value = -31200
$ws.Range("H67").Value = 30750
$ws.Range("J67").Value = 30750
$ws.Range("L67").Value = 30750
$ws.Range("N67").Value = -32310
$ws.Range("H103").Value = 26937.5
$ws.Range("J103").Value = 26937.5
$ws.Range("L103").Value = 26937.5
$ws.Range("N103").Value = -29281.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 30000
$ws.Range("J82").Value = 30000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30766
$ws.Range("H85").Value = 30000
$ws.Range("J85").Value = 30000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32652
$ws.Range("H92").Value = 29850
$ws.Range("J92").Value = 29850
$ws.Range("L92").Value = 29850
$ws.Range("N92").Value = -34842
$ws.Range("H100").Value = 942.65216
$ws.Range("I100").Value = 949.1429000000001
$ws.Range("J100").Value = 932.55554
$ws.Range("K100").Value = 1898.2858
$ws.Range("L100").Value = 1865.11108
$ws.Range("M100").Value = -1357.2858
$ws.Range("N100").Value = -2947.11108
$ws.Range("H101").Value = 22301
$ws.Range("J101").Value = 22301
$ws.Range("L101").Value = 22301
$ws.Range("N101").Value = -28791
$ws.Range("H103").Value = 29879.143
$ws.Range("J103").Value = 29879.143
$ws.Range("L103").Value = 29879.143
$ws.Range("N103").Value = -32223.143
